# Updated cryptos list — apply new Price (D) and Volume(1h) (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so decimal-looking strings like "19.70"
# or "0.0810" keep their exact text (incl. trailing zeros) instead of being
# auto-coerced to numbers by Excel's normal cell-entry parsing.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "49.043.70"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").Value = "2.627.42"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "111.35"
$ws.Range("E5").Value = "  +1.17%  "
$ws.Range("D6").Value = "322.38"
$ws.Range("E6").Value = "  -2.51%  "
$ws.Range("D7").Value = "0.526"
$ws.Range("E7").Value = "  -1.73%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "0.541"
$ws.Range("E9").Value = "  -3.40%  "
$ws.Range("D10").Value = "39.64"
$ws.Range("E10").Value = "  -2.97%  "
$ws.Range("D11").Value = "19.70"
$ws.Range("E11").Value = "  -5.12%  "
$ws.Range("D12").Value = "0.0810"
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").Value = "7.24"
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").Value = "3.035.96"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").Value = "2.634.60"
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("D17").Value = "0.858"
$ws.Range("E17").Value = "  -2.13%  "
$ws.Range("D18").Value = "49.011.37"
$ws.Range("E18").Value = "  -1.40%  "
$ws.Range("D19").Value = "3.02"
$ws.Range("E19").Value = "  -3.04%  "
$ws.Range("D20").Value = "12.91"
$ws.Range("E20").Value = "  -3.53%  "
$ws.Range("E21").Value = "  -2.06%  "
$ws.Range("E22").Value = "  -1.30%  "
$ws.Range("D23").Value = "269.06"
$ws.Range("E23").Value = "  -4.53%  "
$ws.Range("D24").Value = "68.51"
$ws.Range("E24").Value = "  -5.91%  "
$ws.Range("D25").Value = "2.53"
$ws.Range("E25").Value = "  -2.85%  "
$ws.Range("D26").Value = "26.12"
$ws.Range("E26").Value = "  -2.02%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").Value = "10.17"
$ws.Range("E28").Value = "  +1.89%  "
$ws.Range("E29").Value = "  -0.58%  "
$ws.Range("D30").Value = "35.09"
$ws.Range("E30").Value = "  -3.22%  "
$ws.Range("E31").Value = "  -4.15%  "
$ws.Range("D32").Value = "49.40"
$ws.Range("E32").Value = "  -0.79%  "
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("D35").Value = "0.0800"
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("D36").Value = "19.02"
$ws.Range("E36").Value = "  -3.85%  "
$ws.Range("E37").Value = "  +4.23%  "
$ws.Range("D38").Value = "2.04"
$ws.Range("E38").Value = "  -1.19%  "
$ws.Range("D39").Value = "3.12"
$ws.Range("E39").Value = "  +0.96%  "
$ws.Range("D40").Value = "126.82"
$ws.Range("E40").Value = "  +2.80%  "
$ws.Range("E41").Value = "  -1.83%  "
$ws.Range("D42").Value = "22.14"
$ws.Range("E42").Value = "  -3.11%  "
$ws.Range("E43").Value = "  -4.24%  "
$ws.Range("D44").Value = "0.0317"
$ws.Range("E44").Value = "  +0.37%  "
$ws.Range("D45").Value = "2.064.36"
$ws.Range("E45").Value = "  +0.59%  "
$ws.Range("D46").Value = "2.15"
$ws.Range("E46").Value = "  +6.88%  "
$ws.Range("E47").Value = "  -5.02%  "
$ws.Range("E48").Value = "  -5.28%  "
$ws.Range("E49").Value = "  -1.67%  "
$ws.Range("D50").Value = "5.19"
$ws.Range("E50").Value = "  -3.50%  "
$ws.Range("D51").Value = "58.56"
$ws.Range("E51").Value = "  +1.21%  "

# Restore the default cell style on column D (matches original workbook,
# which used the default "Normal" style with no explicit number format).
$dRange.Style = "Normal"

